$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark (it will be recreated at the new
# last-edit position, mirroring what Word does automatically).
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Find the "Programmierer" run and position right after it (before the
# following ":" run) to insert the new "/UX-Manager" text.
$rng = $d.Content
$found = $rng.Find.Execute("Programmierer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPoint = $rng.Duplicate
    $insertPoint.Collapse(0)
    $insertPoint.InsertAfter("/UX-Manager")

    # New _GoBack bookmark marks the position right after the inserted text.
    $bmRange = $insertPoint.Duplicate
    $bmRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
